$wb = $excel.ActiveWorkbook

# --- "Lookup Table" sheet: remove the blank separator row (old row 40). ---
# This shifts the old rows 41-87 up to become the new rows 40-86, which
# also updates the sheet dimension from A1:G87 to A1:G86 automatically.
$ws = $wb.Worksheets.Item("Lookup Table")
$ws.Rows("40").Delete()

# After the shift, two cells that were previously blank (old rows 65 & 66,
# now rows 64 & 65) need their "clmn_num" value (column C) populated.
$ws.Range("C64").Value = "0100"
$ws.Range("C65").Value = "0200"

# --- "Type and Label" sheet: fill in the "type" column (B) for rows 26-30
# with "alpha", matching the other stock/flow rows around them. ---
$ws2 = $wb.Worksheets.Item("Type and Label")
$ws2.Range("B26").Value = "alpha"
$ws2.Range("B27").Value = "alpha"
$ws2.Range("B28").Value = "alpha"
$ws2.Range("B29").Value = "alpha"
$ws2.Range("B30").Value = "alpha"
